# Swap the first two comma-separated names/emails in column G ("Recorded By")
# for every data row in the active worksheet. Rows whose value contains only
# a single entry (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the 7th column ("Recorded By")
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
